$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text so numeric-looking strings (e.g. "1.00", "8.90")
# are preserved exactly instead of being parsed into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '96.540.74'
$ws.Range("E2").Value = '  -0.61%  '
$ws.Range("D3").Value = '3.667.66'
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = '242.65'
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("D6").Value = '1.88'
$ws.Range("E6").Value = '  +13.33%  '
$ws.Range("D7").Value = '663.12'
$ws.Range("E7").Value = '  +1.23%  '
$ws.Range("D8").Value = '0.425'
$ws.Range("E8").Value = '  +3.50%  '
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("E10").Value = '  -0.04%  '
$ws.Range("D11").Value = '3.662.53'
$ws.Range("E11").Value = '  +2.27%  '
$ws.Range("D12").Value = '45.37'
$ws.Range("E12").Value = '  +4.29%  '
$ws.Range("D13").Value = '0.205'
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").Value = '6.86'
$ws.Range("E14").Value = '  +6.20%  '
$ws.Range("D15").Value = '4.347.29'
$ws.Range("E15").Value = '  +2.18%  '
$ws.Range("D16").Value = '0.0000271'
$ws.Range("E16").Value = '  +5.58%  '
$ws.Range("D17").Value = '96.305.03'
$ws.Range("E17").Value = '  -0.61%  '
$ws.Range("D18").Value = '8.97'
$ws.Range("E18").Value = '  +15.55%  '
$ws.Range("D19").Value = '3.657.29'
$ws.Range("E19").Value = '  +1.94%  '
$ws.Range("D20").Value = '12.84'
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").Value = '18.37'
$ws.Range("E21").Value = '  +1.89%  '
$ws.Range("D22").Value = '0.533'
$ws.Range("E22").Value = '  +0.87%  '
$ws.Range("D23").Value = '524.17'
$ws.Range("E23").Value = '  +2.74%  '
$ws.Range("D24").Value = '3.45'
$ws.Range("E24").Value = '  +1.24%  '
$ws.Range("D25").Value = '0.0000205'
$ws.Range("E25").Value = '  +2.16%  '
$ws.Range("D26").Value = '6.99'
$ws.Range("E26").Value = '  +2.07%  '
$ws.Range("D27").Value = '102.01'
$ws.Range("E27").Value = '  +5.13%  '
$ws.Range("D28").Value = '13.10'
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.862.55'
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '0.167'
$ws.Range("E30").Value = '  +10.84%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '12.53'
$ws.Range("E31").Value = '  +8.81%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '3.05'
$ws.Range("E32").Value = '  +0.66%  '
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").Value = '0.187'
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("D35").Value = '1.84'
$ws.Range("E35").Value = '  +13.76%  '
$ws.Range("D36").Value = '33.02'
$ws.Range("E36").Value = '  +4.44%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = '0.591'
$ws.Range("E38").Value = '  +3.82%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '632.25'
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = '8.82'
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '44.93'
$ws.Range("E41").Value = '  +35.19%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.161'
$ws.Range("E42").Value = '  +5.46%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '0.965'
$ws.Range("E43").Value = '  +6.03%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = '1.97'
$ws.Range("E44").Value = '  +5.32%  '
$ws.Range("D45").Value = '6.40'
$ws.Range("E45").Value = '  +9.16%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '0.0460'
$ws.Range("E47").Value = '  +6.49%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.446'
$ws.Range("E48").Value = '  +18.12%  '
$ws.Range("D49").Value = '2.30'
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").Value = '23.63'
$ws.Range("E50").Value = '  -0.20%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '8.62'
$ws.Range("E51").Value = '  +3.30%  '

# Restore the original (default) cell style now that the text values are locked in,
# so no stray style index is left on the data cells.
$ws.Range("D2:E51").Style = "Normal"

